$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Given Input" value in row 5 (special characters test case)
$ws.Range("C5").Value = "$ ass 1.2"

# Update "Expected Output" and "Actual Output" values in row 5
$ws.Range("D5").Value = "char 1 string 3 float 8"
$ws.Range("E5").Value = "char 1 string 3 float 8"

# Update the active cell selection shown in the sheet view
$ws.Range("D7").Select()
